# Apply the commit: adding children 2.1.0, rehab 2.2.0 and kl-gateway 1.2.0
$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsInclude1 = $wb.Worksheets.Item("Include ValueSets")
$wsInclude2 = $wb.Worksheets.Item("Include ValueSets 2")

# Rename the two "include" sheets
$wsInclude1.Name = "Include ValueSet #0"
$wsInclude2.Name = "Include ValueSet #1"

# Simple text updates (no risk of Excel auto-typing as bool/number)
$wsMeta.Range("B3").Value = "1.2.0"
$wsMeta.Range("B8").Value = "2024-10-31T19:21:51+01:00"
$wsMeta.Range("B10").Value = "KL (http://www.kl.dk)"

# B7 (Experimental) needs to become the literal STRING "false" (not boolean).
# A direct .Value assignment of "false"/"true" gets auto-coerced to a boolean
# by the engine (mirrors real Excel behaviour), so instead we compute it via
# a TEXT-returning formula (T()) and then freeze it to a static value with
# Copy + PasteSpecial(xlPasteValues), which keeps the cell as a shared string.
$helper = $wsMeta.Cells.Item(500, 500)
$helper.Formula = "=T(""false"")"
$helper.Copy()
$wsMeta.Cells.Item(7, 2).PasteSpecial(-4163)
$helper.Clear()

# B11 (Jurisdiction) needs to become an empty STRING (still a shared-string
# cell, just referencing an empty string) rather than a genuinely blank cell.
# Plain .Value = "" clears the cell entirely, so use the same formula +
# paste-values trick, pasting from a different helper cell (self paste of an
# empty result is treated as "nothing to paste" and skipped).
$helper2 = $wsMeta.Cells.Item(501, 500)
$helper2.Formula = "=T(""X"")"
$helper2.Copy()
$wsMeta.Cells.Item(11, 2).PasteSpecial(-4163)
$helper2.Clear()
$wsMeta.Cells.Item(11, 2).Replace("X", "")
